# Updates crypto Price (D) / Volume(1h) (E) columns for the refreshed snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.974.07"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "1.677.48"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'215.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "  +1.42%  "
$ws.Range("E8").Value = "  -0.33%  "
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("D10").Value = "'20.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.76%  "
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").Value = "1.915.47"
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("D13").Value = "1.692.58"
$ws.Range("E13").Value = "  +1.26%  "
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("D16").Value = "'65.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "26.988.71"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("D18").Value = "'237.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("E19").Value = "  +4.27%  "
$ws.Range("D20").Value = "0.0₃0734"
$ws.Range("E20").Value = "  -0.54%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  -0.70%  "
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("E24").Value = "  -1.56%  "
$ws.Range("D25").Value = "'145.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("D26").Value = "'7.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.32%  "
$ws.Range("E27").Value = "  +0.71%  "
$ws.Range("D28").Value = "'0.113"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.25%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("D33").Value = "1.488.29"
$ws.Range("E33").Value = "  +1.19%  "
$ws.Range("E34").Value = "  +0.80%  "
$ws.Range("E35").Value = "  +4.50%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("D37").Value = "'0.588"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.36%  "
$ws.Range("D38").Value = "'0.0174"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.80%  "
$ws.Range("D39").Value = "'0.902"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.43%  "
$ws.Range("E40").Value = "  -4.19%  "
$ws.Range("E41").Value = "  +0.60%  "
$ws.Range("E43").Value = "  +1.89%  "
$ws.Range("D44").Value = "'67.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.50%  "
$ws.Range("D45").Value = "1.820.63"
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("D46").Value = "'0.782"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.23%  "
$ws.Range("D47").Value = "'90.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.43%  "
$ws.Range("D48").Value = "0.0₆0106"
$ws.Range("E48").Value = "  +13.10%  "
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("E50").Value = "  +1.85%  "
$ws.Range("E51").Value = "  +0.40%  "
